# Apply a growth constraint for wave and tidal (PWR sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PWR")

# --- Row 6 (CCS): switch the label/description cells to formulas driven by A12 ---
$ws.Range("B6").Formula = '=TEXTJOIN("",TRUE,"UC-PWR_MaxGrowth",A12)'
$ws.Range("J6").Formula = '=1+$B12'
$ws.Range("N6").Formula = '=TEXTJOIN(" ",TRUE,"PWR",A12, "maximum growth rate")'

# --- Row 7 (Wave) ---
$ws.Range("B7").Formula = '=TEXTJOIN("",TRUE,"UC-PWR_MaxGrowth",A13)'
$ws.Range("C7").Value = "CAP, GROWTH"
$ws.Range("G7").Value = "P*OCE*TID*"
$ws.Range("H7").Value = 2035
$ws.Range("I7").Value = "LO"
$ws.Range("J7").Formula = '=1+$B13'
$ws.Range("K7").Value = 1
$ws.Range("L7").Formula = '=-D13'
$ws.Range("M7").Value = 5
$ws.Range("N7").Formula = '=TEXTJOIN(" ",TRUE,"PWR",A13, "maximum growth rate")'

# --- Row 8 (Tidal) ---
$ws.Range("B8").Formula = '=TEXTJOIN("",TRUE,"UC-PWR_MaxGrowth",A14)'
$ws.Range("C8").Value = "CAP, GROWTH"
$ws.Range("G8").Value = "P*OCE*WAV*"
$ws.Range("H8").Value = 2035
$ws.Range("I8").Value = "LO"
$ws.Range("J8").Formula = '=1+$B14'
$ws.Range("K8").Value = 1
$ws.Range("L8").Formula = '=-D14'
$ws.Range("M8").Value = 5
$ws.Range("N8").Formula = '=TEXTJOIN(" ",TRUE,"PWR",A14, "maximum growth rate")'

# --- Row 12 (CCS starting data): add the set label in column A ---
$ws.Range("A12").Value = "CCS"

# --- Row 13 (Wave starting data) ---
$ws.Range("A13").Value = "Wave"
$ws.Range("B13").Value = 0.2
$ws.Range("D13").Value = 0.4

# --- Row 14 (Tidal starting data) ---
$ws.Range("A14").Value = "Tidal"
$ws.Range("B14").Value = 0.2
$ws.Range("D14").Value = 0.1

# Copy number formats/styles for the new rows from their CCS-row analogues
$ws.Range("B6:N6").Copy()
$ws.Range("B7:N7").PasteSpecial(-4122)
$ws.Range("B6:N6").Copy()
$ws.Range("B8:N8").PasteSpecial(-4122)
$ws.Range("B12:D12").Copy()
$ws.Range("B13:D13").PasteSpecial(-4122)
$ws.Range("B12:D12").Copy()
$ws.Range("B14:D14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the values/formulas (PasteSpecial of formats only shouldn't disturb these,
# but re-assert to be safe against any paste side-effects)
$ws.Range("B7").Formula = '=TEXTJOIN("",TRUE,"UC-PWR_MaxGrowth",A13)'
$ws.Range("C7").Value = "CAP, GROWTH"
$ws.Range("G7").Value = "P*OCE*TID*"
$ws.Range("H7").Value = 2035
$ws.Range("I7").Value = "LO"
$ws.Range("J7").Formula = '=1+$B13'
$ws.Range("K7").Value = 1
$ws.Range("L7").Formula = '=-D13'
$ws.Range("M7").Value = 5
$ws.Range("N7").Formula = '=TEXTJOIN(" ",TRUE,"PWR",A13, "maximum growth rate")'

$ws.Range("B8").Formula = '=TEXTJOIN("",TRUE,"UC-PWR_MaxGrowth",A14)'
$ws.Range("C8").Value = "CAP, GROWTH"
$ws.Range("G8").Value = "P*OCE*WAV*"
$ws.Range("H8").Value = 2035
$ws.Range("I8").Value = "LO"
$ws.Range("J8").Formula = '=1+$B14'
$ws.Range("K8").Value = 1
$ws.Range("L8").Formula = '=-D14'
$ws.Range("M8").Value = 5
$ws.Range("N8").Formula = '=TEXTJOIN(" ",TRUE,"PWR",A14, "maximum growth rate")'

$ws.Range("A12").Value = "CCS"
$ws.Range("A13").Value = "Wave"
$ws.Range("B13").Value = 0.2
$ws.Range("D13").Value = 0.4
$ws.Range("A14").Value = "Tidal"
$ws.Range("B14").Value = 0.2
$ws.Range("D14").Value = 0.1

# --- Column widths (closest representable values on this engine's pixel grid) ---
$ws.Columns("B").ColumnWidth = 24.666666666666668
$ws.Columns("G").ColumnWidth = 12.5
$ws.Columns("N").ColumnWidth = 30.5

# --- Selection (cosmetic) ---
$ws.Range("N16").Select()
